$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before the old column Y ("Khách đã trả"),
# i.e. right after column X ("Khách hàng cần trả"). This shifts every
# column from Y onward two places to the right and carries over the
# formatting of column Y to the two new columns - matching the
# "BuTruTraHang" column-add commit.
$ws.Columns("Y:Z").Insert()

# New header captions for the two inserted columns.
$ws.Range("Y6").Value2 = "Bù trừ trả hàng"
$ws.Range("Z6").Value2 = "Giá trị sau trả"

# Totals row (row 31) needs the same SUM formulas as the other data
# columns for the two newly inserted columns.
$ws.Range("Y31").Formula = "=SUM(Y`$7:Y30)"
$ws.Range("Z31").Formula = "=SUM(Z`$7:Z30)"

# Match the row height bump on the totals row seen in the target file.
$ws.Rows(31).RowHeight = 24.75

# Restore the active cell in the frozen (bottom-left) pane.
$ws.Range("J12").Select()
